# Add a new "Player Info" sheet before the existing "ODI Batting" sheet,
# populate it with the player's biographical data, and update the
# "ODI Batting" sheet so the MATCH_CARD_LINK column becomes MATCH_CODE
# (storing just the numeric match code instead of the full scorecard URL).

$wb = $excel.ActiveWorkbook
$battingSheet = $wb.Worksheets.Item(1)

# --- 1. Insert the new "Player Info" worksheet as the first tab ---------
$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# ID is stored as text (e.g. "6630"), not a number, so format the cell as
# text first -- otherwise Excel auto-coerces the numeric-looking string.
$infoSheet.Range("A2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "6630"
$infoSheet.Range("B2").Value = "Mohammad Haris"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Off Break"

# Match the bold / bordered / centered header look used on the other sheet.
$headerRange = $infoSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Re-resolve the original sheet by name: inserting a sheet re-numbers the
# 1-based worksheet index, so the old $battingSheet reference now points
# at whatever sheet occupies its former slot (i.e. the newly added one).
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- 2. Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" -----------
$battingSheet.Range("D1").Value = "MATCH_CODE"

# --- 3. Replace full scorecard URLs with the bare numeric match code ----
# These stay text cells (like the rest of the column previously holding
# the URLs), so force the text format before writing the digit strings.
$codeRange = $battingSheet.Range("D2:D5")
$codeRange.NumberFormat = "@"
$battingSheet.Range("D2").Value = "4586"
$battingSheet.Range("D3").Value = "4590"
$battingSheet.Range("D4").Value = "4592"
$battingSheet.Range("D5").Value = "4641"
